# Adds Industry Category subscript to SoCiIEPTtB Share of Change in Industry
# Expenses Passed Through to Buyers; values set to 1 for all industries.

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("SoCiIEPTtB")

# ---------------------------------------------------------------------------
# 1. SoCiIEPTtB data sheet: expand the single "Share of expenses" row into a
#    full Industry Category subscript (25 industries), each with a formula
#    pulling the passthrough share from About!A42 (always 1 / full
#    passthrough).
# ---------------------------------------------------------------------------

$industries = @(
    "agriculture and forestry 01T03",
    "coal mining 05",
    "oil and gas extraction 06",
    "other mining and quarrying 07T08",
    "food beverage and tobacco 10T12",
    "textiles apparel and leather 13T15",
    "wood products 16",
    "pulp paper and printing 17T18",
    "refined petroleum and coke 19",
    "chemicals 20",
    "rubber and plastic products 22",
    "glass and glass products 231",
    "cement and other nonmetallic minerals 239",
    "iron and steel 241",
    "other metals 242",
    "metal products except machinery and vehicles 25",
    "computers and electronics 26",
    "appliances and electrical equipment 27",
    "other machinery 28",
    "road vehicles 29",
    "nonroad vehicles 30",
    "other manufacturing 31T33",
    "energy pipelines and gas processing 352T353",
    "water and waste 36T39",
    "construction 41T43"
)

# Headers stay the same text; only the data rows below them change.
$wsData.Range("A1").Value = "Unit: dimensionless (% passthrough)"
$wsData.Range("B1").Value = "Paid by consumers"

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 1).Value = $industries[$i]
    $wsData.Cells.Item($row, 2).Formula = "=About!A`$42"
}

# ---------------------------------------------------------------------------
# 2. About sheet: the two footnote-style rows at the bottom (A40/A41) lose
#    their (unused/empty-alignment) cell style.
# ---------------------------------------------------------------------------

$wsAbout.Range("A40:A41").ClearFormats()

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping: About!B42 is selected (but About
#    is no longer the active tab), and SoCiIEPTtB becomes the active tab with
#    I16 selected.
# ---------------------------------------------------------------------------

$wsAbout.Activate()
$wsAbout.Range("B42").Select()

$wsData.Activate()
$wsData.Range("I16").Select()
